$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H12").Value = 761.3333
$ws.Range("I12").Value = 364.14285
$ws.Range("J12").Value = 1317.4
$ws.Range("K12").Value = 364.14285
$ws.Range("L12").Value = 1317.4
$ws.Range("M12").Value = -194.14285
$ws.Range("N12").Value = -1657.4
$ws.Range("H31").Value = 609.625
$ws.Range("I31").Value = 479.5
$ws.Range("K31").Value = 1438.5
$ws.Range("M31").Value = -1208.5
$ws.Range("H55").Value = 544.4545000000001
$ws.Range("I55").Value = 385.66666
$ws.Range("J55").Value = 654.38464
$ws.Range("K55").Value = 385.66666
$ws.Range("L55").Value = 654.38464
$ws.Range("M55").Value = -171.66666
$ws.Range("N55").Value = -1082.38464
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H113").Value = 75867.86
$ws.Range("I113").Value = 251787.75
$ws.Range("J113").Value = 5499.9
$ws.Range("K113").Value = 251787.75
$ws.Range("L113").Value = 5499.9
$ws.Range("M113").Value = -248533.75
$ws.Range("N113").Value = -12007.9
$ws.Range("H116").Value = 8349.237999999999
$ws.Range("I116").Value = 11169.583
$ws.Range("K116").Value = 11169.583
$ws.Range("M116").Value = -7727.583000000001
$ws.Range("H137").Value = 3092
$ws.Range("J137").Value = 3092
$ws.Range("L137").Value = 9276
$ws.Range("N137").Value = -14376

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1685.8649
$ws.Range("I2").Value = 1645.6061
$ws.Range("K2").Value = 1645.6061
$ws.Range("M2").Value = -1532.6061
$ws.Range("H32").Value = 4572.2983
$ws.Range("I32").Value = 3862.635
$ws.Range("K32").Value = 3862.635
$ws.Range("M32").Value = -3575.635
$ws.Range("H61").Value = 1749.9412
$ws.Range("I61").Value = 1749.9412
$ws.Range("K61").Value = 1749.9412
$ws.Range("M61").Value = -1537.9412
$ws.Range("H63").Value = 5584.8335
$ws.Range("I63").Value = 3877.25
$ws.Range("J63").Value = 9000
$ws.Range("K63").Value = 3877.25
$ws.Range("L63").Value = 9000
$ws.Range("M63").Value = -3191.25
$ws.Range("N63").Value = -10372
$ws.Range("H66").Value = 5584.8335
$ws.Range("I66").Value = 3877.25
$ws.Range("J66").Value = 9000
$ws.Range("K66").Value = 19386.25
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -15954.25
$ws.Range("N66").Value = -51864
$ws.Range("H102").Value = 1557.0968
$ws.Range("I102").Value = 1652.1072
$ws.Range("J102").Value = 670.3333
$ws.Range("K102").Value = 1652.1072
$ws.Range("L102").Value = 670.3333
$ws.Range("M102").Value = -30.10719999999992
$ws.Range("N102").Value = -3914.3333
$ws.Range("H116").Value = 1685.8649
$ws.Range("I116").Value = 1645.6061
$ws.Range("K116").Value = 1645.6061
$ws.Range("M116").Value = 648.3939
$ws.Range("H122").Value = 1641.48
$ws.Range("I122").Value = 1222.6487
$ws.Range("J122").Value = 2833.5386
$ws.Range("K122").Value = 3667.9461
$ws.Range("L122").Value = 8500.6158
$ws.Range("M122").Value = -1217.9461
$ws.Range("N122").Value = -13400.6158
$ws.Range("H132").Value = 5445.8667
$ws.Range("I132").Value = 4844.3774
$ws.Range("K132").Value = 14533.1322
$ws.Range("M132").Value = -12003.1322
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 1749.9412
$ws.Range("I136").Value = 1749.9412
$ws.Range("K136").Value = 5249.8236
$ws.Range("M136").Value = -2699.8236

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1685.8649
$ws.Range("I3").Value = 1645.6061
$ws.Range("K3").Value = 1645.6061
$ws.Range("M3").Value = -1531.6061
$ws.Range("H86").Value = 3794
$ws.Range("I86").Value = 3546.5334
$ws.Range("J86").Value = 5650
$ws.Range("K86").Value = 3546.5334
$ws.Range("L86").Value = 5650
$ws.Range("M86").Value = -2423.5334
$ws.Range("N86").Value = -7896
$ws.Range("H89").Value = 3794
$ws.Range("I89").Value = 3546.5334
$ws.Range("J89").Value = 5650
$ws.Range("K89").Value = 17732.667
$ws.Range("L89").Value = 28250
$ws.Range("M89").Value = -12116.667
$ws.Range("N89").Value = -39482
$ws.Range("H107").Value = 101683.8
$ws.Range("I107").Value = 143985.58
$ws.Range("K107").Value = 143985.58
$ws.Range("M107").Value = -142065.58
$ws.Range("H132").Value = 130524
$ws.Range("J132").Value = 130524
$ws.Range("L132").Value = 130524
$ws.Range("N132").Value = -140644
$ws.Range("H134").Value = 1169.8125
$ws.Range("I134").Value = 1014.4667
$ws.Range("K134").Value = 3043.4001
$ws.Range("M134").Value = -508.4000999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27638.7
$ws.Range("I31").Value = 29528.75
$ws.Range("J31").Value = 10628.25
$ws.Range("K31").Value = 29528.75
$ws.Range("L31").Value = 10628.25
$ws.Range("M31").Value = -29233.75
$ws.Range("N31").Value = -11218.25
$ws.Range("H34").Value = 27638.7
$ws.Range("I34").Value = 29528.75
$ws.Range("J34").Value = 10628.25
$ws.Range("K34").Value = 29528.75
$ws.Range("L34").Value = 10628.25
$ws.Range("M34").Value = -29326.75
$ws.Range("N34").Value = -11032.25
$ws.Range("H141").Value = 234495.4
$ws.Range("J141").Value = 234495.4
$ws.Range("L141").Value = 234495.4
$ws.Range("N141").Value = -244855.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 589.3
$ws.Range("I2").Value = 85.84614999999999
$ws.Range("J2").Value = 974.2941
$ws.Range("K2").Value = 515.0769
$ws.Range("L2").Value = 5845.7646
$ws.Range("M2").Value = -402.0769
$ws.Range("N2").Value = -6071.7646
$ws.Range("H23").Value = 178.16667
$ws.Range("J23").Value = 144.5
$ws.Range("L23").Value = 433.5
$ws.Range("N23").Value = -903.5
$ws.Range("H37").Value = 333358340
$ws.Range("J37").Value = 333358340
$ws.Range("L37").Value = 1000075020
$ws.Range("N37").Value = -1000075244
$ws.Range("H40").Value = 368
$ws.Range("J40").Value = 366.5
$ws.Range("L40").Value = 1466
$ws.Range("N40").Value = -1604
$ws.Range("H70").Value = 4211.5713
$ws.Range("I70").Value = 2497.25
$ws.Range("J70").Value = 6497.3335
$ws.Range("K70").Value = 7491.75
$ws.Range("L70").Value = 19492.0005
$ws.Range("M70").Value = -7176.75
$ws.Range("N70").Value = -20122.0005
$ws.Range("H73").Value = 4211.5713
$ws.Range("I73").Value = 2497.25
$ws.Range("J73").Value = 6497.3335
$ws.Range("K73").Value = 7491.75
$ws.Range("L73").Value = 19492.0005
$ws.Range("M73").Value = -6399.75
$ws.Range("N73").Value = -21676.0005
$ws.Range("H82").Value = 7257
$ws.Range("I82").Value = 3385.5
$ws.Range("J82").Value = 15000
$ws.Range("K82").Value = 10156.5
$ws.Range("L82").Value = 45000
$ws.Range("M82").Value = -9750.5
$ws.Range("N82").Value = -45812
$ws.Range("H85").Value = 7257
$ws.Range("I85").Value = 3385.5
$ws.Range("J85").Value = 15000
$ws.Range("K85").Value = 10156.5
$ws.Range("L85").Value = 45000
$ws.Range("M85").Value = -8752.5
$ws.Range("N85").Value = -47808
$ws.Range("H113").Value = 2807.1538
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2807.1538
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 8421.4614
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -12761.4614
$ws.Range("H140").Value = 1968.1111
$ws.Range("I140").Value = 1919.1765
$ws.Range("K140").Value = 5757.529500000001
$ws.Range("M140").Value = -577.5295000000006

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 23812152
$ws.Range("I107").Value = 644
$ws.Range("J107").Value = 55560828
$ws.Range("K107").Value = 644
$ws.Range("L107").Value = 55560828
$ws.Range("M107").Value = 1276
$ws.Range("N107").Value = -55564668
$ws.Range("H132").Value = 2322.6943
$ws.Range("I132").Value = 2322.6943
$ws.Range("K132").Value = 6968.0829
$ws.Range("M132").Value = -4438.0829

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 125804.125
$ws.Range("I55").Value = 167638.83
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 167638.83
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = -167465.83
$ws.Range("N55").Value = -646
$ws.Range("H95").Value = 17499.5
$ws.Range("J95").Value = 17499.5
$ws.Range("L95").Value = 17499.5
$ws.Range("N95").Value = -22991.5
$ws.Range("H109").Value = 21070.715
$ws.Range("J109").Value = 21070.715
$ws.Range("L109").Value = 21070.715
$ws.Range("N109").Value = -23844.715
$ws.Range("H122").Value = 87229.336
$ws.Range("J122").Value = 6190
$ws.Range("L122").Value = 18570
$ws.Range("N122").Value = -23470
$ws.Range("H132").Value = 2976.5334
$ws.Range("I132").Value = 2620.5435
$ws.Range("K132").Value = 7861.630500000001
$ws.Range("M132").Value = -5331.630500000001
$ws.Range("H136").Value = 5735
$ws.Range("I136").Value = 4981.6665
$ws.Range("K136").Value = 14944.9995
$ws.Range("M136").Value = -12394.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 38333
$ws.Range("J97").Value = 38333
$ws.Range("L97").Value = 38333
$ws.Range("N97").Value = -40315
$ws.Range("H107").Value = 17247322
$ws.Range("I107").Value = 7060.421
$ws.Range("K107").Value = 21181.263
$ws.Range("M107").Value = -19261.263
$ws.Range("H132").Value = 3110.4614
$ws.Range("I132").Value = 2518.9111
$ws.Range("K132").Value = 7556.7333
$ws.Range("M132").Value = -5026.7333
